$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 (I0) and J1 (IF), copying the header style from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in I and J column data for rows 2-38
$data = @{
    2 = @(7, 8)
    3 = @(9, 9)
    4 = @(10, 10)
    5 = @(7, 7)
    6 = @(7, 8)
    7 = @(6, 7)
    8 = @(5, 5)
    9 = @(9, 9)
    10 = @(6, 6)
    11 = @(8, 8)
    12 = @(9, 9)
    13 = @(8, 8)
    14 = @(7, 8)
    15 = @(7, 7)
    16 = @(7, 7)
    17 = @(8, 8)
    18 = @(8, 8)
    19 = @(6, 7)
    20 = @(8, 8)
    21 = @(7, 7)
    22 = @(6, 6)
    23 = @(7, 7)
    24 = @(8, 8)
    25 = @(7, 8)
    26 = @(6, 6)
    27 = @(7, 7)
    28 = @(8, 8)
    29 = @(7, 7)
    30 = @(8, 8)
    31 = @(7, 7)
    32 = @(5, 5)
    33 = @(6, 6)
    34 = @(4, 5)
    35 = @(8, 8)
    36 = @(7, 7)
    37 = @(7, 7)
    38 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}

Write-Host "done"
